# Vtn-Itgb8.xlsx -- refresh LR-pair (NATMI) table with updated TPM-derived values
# and one additional Sending-cluster/Target-cluster combination (Resolving-Mac).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = 17
$cols = 20
$data = New-Object 'object[,]' $rows,$cols

# row 1
$data[0,0] = 'Sending cluster'
$data[0,1] = 'Ligand symbol'
$data[0,2] = 'Receptor symbol'
$data[0,3] = 'Target cluster'
$data[0,4] = 'Ligand-expressing cells'
$data[0,5] = 'Ligand detection rate'
$data[0,6] = 'Ligand average expression value'
$data[0,7] = 'Ligand total expression value'
$data[0,8] = 'Ligand derived specificity of average expression value'
$data[0,9] = 'Ligand derived specificity of total expression value'
$data[0,10] = 'Receptor-expressing cells'
$data[0,11] = 'Receptor detection rate'
$data[0,12] = 'Receptor average expression value'
$data[0,13] = 'Receptor total expression value'
$data[0,14] = 'Receptor derived specificity of average expression value'
$data[0,15] = 'Receptor derived specificity of total expression value'
$data[0,16] = 'Edge average expression weight'
$data[0,17] = 'Edge total expression weight'
$data[0,18] = 'Edge average expression derived specificity'
$data[0,19] = 'Edge total expression derived specificity'

# row 2
$data[1,0] = 'ECs'
$data[1,1] = 'Vtn'
$data[1,2] = 'Itgb8'
$data[1,3] = 'ECs'
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = 3.337313
$data[1,7] = 10.011939
$data[1,8] = 0.1958858017947999
$data[1,9] = 0.1958858017947999
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 0.111675
$data[1,13] = 0.335025
$data[1,14] = 0.02767755395605
$data[1,15] = 0.02767755395605
$data[1,16] = 0.372694429275
$data[1,17] = 3.354249863475
$data[1,18] = 0.005421639848399691
$data[1,19] = 0.005421639848399691

# row 3
$data[2,0] = 'ECs'
$data[2,1] = 'Vtn'
$data[2,2] = 'Itgb8'
$data[2,3] = 'FAPs'
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 3.337313
$data[2,7] = 10.011939
$data[2,8] = 0.1958858017947999
$data[2,9] = 0.1958858017947999
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 2.253188666666667
$data[2,13] = 6.759566
$data[2,14] = 0.5584307221385899
$data[2,15] = 0.5584307221385899
$data[2,16] = 7.519595828719333
$data[2,17] = 67.67636245847399
$data[2,18] = 0.1093886497529668
$data[2,19] = 0.1093886497529668

# row 4
$data[3,0] = 'ECs'
$data[3,1] = 'Vtn'
$data[3,2] = 'Itgb8'
$data[3,3] = 'Inflammatory-Mac'
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 3.337313
$data[3,7] = 10.011939
$data[3,8] = 0.1958858017947999
$data[3,9] = 0.1958858017947999
$data[3,10] = 1
$data[3,11] = 0.3333333333333333
$data[3,12] = 0.01559133333333333
$data[3,13] = 0.046774
$data[3,14] = 0.003864159118693479
$data[3,15] = 0.003864159118693479
$data[3,16] = 0.05203315942066667
$data[3,17] = 0.468298434786
$data[3,18] = 0.0007569339072279596
$data[3,19] = 0.0007569339072279596

# row 5
$data[4,0] = 'ECs'
$data[4,1] = 'Vtn'
$data[4,2] = 'Itgb8'
$data[4,3] = 'MuSCs'
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 3.337313
$data[4,7] = 10.011939
$data[4,8] = 0.1958858017947999
$data[4,9] = 0.1958858017947999
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 1.654403
$data[4,13] = 4.963209
$data[4,14] = 0.4100275647866666
$data[4,15] = 0.4100275647866666
$data[4,16] = 5.521260639139
$data[4,17] = 49.691345752251
$data[4,18] = 0.08031857828620546
$data[4,19] = 0.08031857828620546

# row 6
$data[5,0] = 'FAPs'
$data[5,1] = 'Vtn'
$data[5,2] = 'Itgb8'
$data[5,3] = 'ECs'
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 10.28369633333333
$data[5,7] = 30.851089
$data[5,8] = 0.6036083824529627
$data[5,9] = 0.6036083824529627
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 0.111675
$data[5,13] = 0.335025
$data[5,14] = 0.02767755395605
$data[5,15] = 0.02767755395605
$data[5,16] = 1.148431788025
$data[5,17] = 10.335886092225
$data[5,18] = 0.01670640357366594
$data[5,19] = 0.01670640357366594

# row 7
$data[6,0] = 'FAPs'
$data[6,1] = 'Vtn'
$data[6,2] = 'Itgb8'
$data[6,3] = 'FAPs'
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 10.28369633333333
$data[6,7] = 30.851089
$data[6,8] = 0.6036083824529627
$data[6,9] = 0.6036083824529627
$data[6,10] = 3
$data[6,11] = 1
$data[6,12] = 2.253188666666667
$data[6,13] = 6.759566
$data[6,14] = 0.5584307221385899
$data[6,15] = 0.5584307221385899
$data[6,16] = 23.17110802970822
$data[6,17] = 208.539972267374
$data[6,18] = 0.3370734649021141
$data[6,19] = 0.3370734649021141

# row 8
$data[7,0] = 'FAPs'
$data[7,1] = 'Vtn'
$data[7,2] = 'Itgb8'
$data[7,3] = 'Inflammatory-Mac'
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 10.28369633333333
$data[7,7] = 30.851089
$data[7,8] = 0.6036083824529627
$data[7,9] = 0.6036083824529627
$data[7,10] = 1
$data[7,11] = 0.3333333333333333
$data[7,12] = 0.01559133333333333
$data[7,13] = 0.046774
$data[7,14] = 0.003864159118693479
$data[7,15] = 0.003864159118693479
$data[7,16] = 0.1603365374317778
$data[7,17] = 1.443028836886
$data[7,18] = 0.002332438835175437
$data[7,19] = 0.002332438835175437

# row 9
$data[8,0] = 'FAPs'
$data[8,1] = 'Vtn'
$data[8,2] = 'Itgb8'
$data[8,3] = 'MuSCs'
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = 10.28369633333333
$data[8,7] = 30.851089
$data[8,8] = 0.6036083824529627
$data[8,9] = 0.6036083824529627
$data[8,10] = 3
$data[8,11] = 1
$data[8,12] = 1.654403
$data[8,13] = 4.963209
$data[8,14] = 0.4100275647866666
$data[8,15] = 0.4100275647866666
$data[8,16] = 17.01337806495567
$data[8,17] = 153.120402584601
$data[8,18] = 0.2474960751420072
$data[8,19] = 0.2474960751420072

# row 10
$data[9,0] = 'MuSCs'
$data[9,1] = 'Vtn'
$data[9,2] = 'Itgb8'
$data[9,3] = 'ECs'
$data[9,4] = 3
$data[9,5] = 1
$data[9,6] = 3.355061
$data[9,7] = 10.065183
$data[9,8] = 0.196927532435664
$data[9,9] = 0.196927532435664
$data[9,10] = 3
$data[9,11] = 1
$data[9,12] = 0.111675
$data[9,13] = 0.335025
$data[9,14] = 0.02767755395605
$data[9,15] = 0.02767755395605
$data[9,16] = 0.3746764371749999
$data[9,17] = 3.372087934574999
$data[9,18] = 0.005450472404419877
$data[9,19] = 0.005450472404419877

# row 11
$data[10,0] = 'MuSCs'
$data[10,1] = 'Vtn'
$data[10,2] = 'Itgb8'
$data[10,3] = 'FAPs'
$data[10,4] = 3
$data[10,5] = 1
$data[10,6] = 3.355061
$data[10,7] = 10.065183
$data[10,8] = 0.196927532435664
$data[10,9] = 0.196927532435664
$data[10,10] = 3
$data[10,11] = 1
$data[10,12] = 2.253188666666667
$data[10,13] = 6.759566
$data[10,14] = 0.5584307221385899
$data[10,15] = 0.5584307221385899
$data[10,16] = 7.559585421175332
$data[10,17] = 68.03626879057799
$data[10,18] = 0.1099703841470184
$data[10,19] = 0.1099703841470184

# row 12
$data[11,0] = 'MuSCs'
$data[11,1] = 'Vtn'
$data[11,2] = 'Itgb8'
$data[11,3] = 'Inflammatory-Mac'
$data[11,4] = 3
$data[11,5] = 1
$data[11,6] = 3.355061
$data[11,7] = 10.065183
$data[11,8] = 0.196927532435664
$data[11,9] = 0.196927532435664
$data[11,10] = 1
$data[11,11] = 0.3333333333333333
$data[11,12] = 0.01559133333333333
$data[11,13] = 0.046774
$data[11,14] = 0.003864159118693479
$data[11,15] = 0.003864159118693479
$data[11,16] = 0.05230987440466666
$data[11,17] = 0.470788869642
$data[11,18] = 0.0007609593201830771
$data[11,19] = 0.0007609593201830771

# row 13
$data[12,0] = 'MuSCs'
$data[12,1] = 'Vtn'
$data[12,2] = 'Itgb8'
$data[12,3] = 'MuSCs'
$data[12,4] = 3
$data[12,5] = 1
$data[12,6] = 3.355061
$data[12,7] = 10.065183
$data[12,8] = 0.196927532435664
$data[12,9] = 0.196927532435664
$data[12,10] = 3
$data[12,11] = 1
$data[12,12] = 1.654403
$data[12,13] = 4.963209
$data[12,14] = 0.4100275647866666
$data[12,15] = 0.4100275647866666
$data[12,16] = 5.550622983583
$data[12,17] = 49.95560685224699
$data[12,18] = 0.08074571656404263
$data[12,19] = 0.08074571656404263

# row 14
$data[13,0] = 'Resolving-Mac'
$data[13,1] = 'Vtn'
$data[13,2] = 'Itgb8'
$data[13,3] = 'ECs'
$data[13,4] = 1
$data[13,5] = 0.3333333333333333
$data[13,6] = 0.06096333333333333
$data[13,7] = 0.18289
$data[13,8] = 0.003578283316573439
$data[13,9] = 0.003578283316573439
$data[13,10] = 3
$data[13,11] = 1
$data[13,12] = 0.111675
$data[13,13] = 0.335025
$data[13,14] = 0.02767755395605
$data[13,15] = 0.02767755395605
$data[13,16] = 0.006808080249999999
$data[13,17] = 0.06127272224999999
$data[13,18] = 0.00009903812956449491
$data[13,19] = 0.0000990381295644949

# row 15
$data[14,0] = 'Resolving-Mac'
$data[14,1] = 'Vtn'
$data[14,2] = 'Itgb8'
$data[14,3] = 'FAPs'
$data[14,4] = 1
$data[14,5] = 0.3333333333333333
$data[14,6] = 0.06096333333333333
$data[14,7] = 0.18289
$data[14,8] = 0.003578283316573439
$data[14,9] = 0.003578283316573439
$data[14,10] = 3
$data[14,11] = 1
$data[14,12] = 2.253188666666667
$data[14,13] = 6.759566
$data[14,14] = 0.5584307221385899
$data[14,15] = 0.5584307221385899
$data[14,16] = 0.1373618917488889
$data[14,17] = 1.23625702574
$data[14,18] = 0.001998223336490574
$data[14,19] = 0.001998223336490574

# row 16
$data[15,0] = 'Resolving-Mac'
$data[15,1] = 'Vtn'
$data[15,2] = 'Itgb8'
$data[15,3] = 'Inflammatory-Mac'
$data[15,4] = 1
$data[15,5] = 0.3333333333333333
$data[15,6] = 0.06096333333333333
$data[15,7] = 0.18289
$data[15,8] = 0.003578283316573439
$data[15,9] = 0.003578283316573439
$data[15,10] = 1
$data[15,11] = 0.3333333333333333
$data[15,12] = 0.01559133333333333
$data[15,13] = 0.046774
$data[15,14] = 0.003864159118693479
$data[15,15] = 0.003864159118693479
$data[15,16] = 0.0009504996511111111
$data[15,17] = 0.00855449686
$data[15,18] = 0.000013827056107006
$data[15,19] = 0.000013827056107006

# row 17
$data[16,0] = 'Resolving-Mac'
$data[16,1] = 'Vtn'
$data[16,2] = 'Itgb8'
$data[16,3] = 'MuSCs'
$data[16,4] = 1
$data[16,5] = 0.3333333333333333
$data[16,6] = 0.06096333333333333
$data[16,7] = 0.18289
$data[16,8] = 0.003578283316573439
$data[16,9] = 0.003578283316573439
$data[16,10] = 3
$data[16,11] = 1
$data[16,12] = 1.654403
$data[16,13] = 4.963209
$data[16,14] = 0.4100275647866666
$data[16,15] = 0.4100275647866666
$data[16,16] = 0.1008579215566667
$data[16,17] = 0.90772129401
$data[16,18] = 0.001467194794411364
$data[16,19] = 0.001467194794411364

$ws.Range("A1:T17").Value2 = $data
